{"js": "// Highlight (green) the five \"Le Dessin\" / \"Le Multijoueur (V1)\" checklist\n// bullet paragraphs that describe the new client-server pixel drawing work:\n//   - \"Cr\u00e9er la palette de couleurs (7 + blanc(gomme))\"\n//   - \"Changer de couleur un pixel\"\n//   - the \"Le Multijoueur (V1)\" sub-heading paragraph (\": \" + bold title)\n//   - \"Connecter le clic au serveur.\"\n//   - \"Tester avec deux onglets ouverts\"\n//\n// Setting `font.highlightColor` on the Paragraph applies the highlight to\n// the paragraph mark (pPr/rPr) as well as every run in the paragraph,\n// matching the OOXML diff (<w:highlight w:val=\"green\"/> added to the\n// paragraph mark rPr and to each run's rPr).\n\nconst targetTexts = [\n  \"Cr\u00e9er la palette de couleurs (7 + blanc(gomme))\",\n  \"Changer de couleur un pixel\",\n  \"Le Multijoueur (V1)\",\n  \"Connecter le clic au serveur.\",\n  \"Tester avec deux onglets ouverts\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  if (targetTexts.some((target) => text.indexOf(target) !== -1)) {\n    paragraph.font.highlightColor = \"#00FF00\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight (green) the five \"Le Dessin\" / \"Le Multijoueur (V1)\" checklist\n# bullet paragraphs that describe the new client-server pixel drawing work:\n#   - \"Cr\u00e9er la palette de couleurs (7 + blanc(gomme))\"\n#   - \"Changer de couleur un pixel\"\n#   - the \"Le Multijoueur (V1)\" sub-heading paragraph (\": \" + bold title)\n#   - \"Connecter le clic au serveur.\"\n#   - \"Tester avec deux onglets ouverts\"\n#\n# $p.Range includes the trailing paragraph mark, so setting\n# Font.HighlightColorIndex on it highlights both the paragraph mark\n# (pPr/rPr) and every run in the paragraph (rPr), matching the OOXML diff\n# (<w:highlight w:val=\"green\"/> added to the paragraph mark rPr and to\n# each run's rPr).\n\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n    \"Cr\u00e9er la palette de couleurs (7 + blanc(gomme))\",\n    \"Changer de couleur un pixel\",\n    \"Le Multijoueur (V1)\",\n    \"Connecter le clic au serveur.\",\n    \"Tester avec deux onglets ouverts\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    foreach ($target in $targetTexts) {\n        if ($text -like \"*$target*\") {\n            $p.Range.Font.HighlightColorIndex = 4\n            break\n        }\n    }\n}\n"}
